# The commit rotates the per-occurrence record data in rows 2-4 down by
# one row (row2 -> row3, row3 -> row4, row4 -> row2); row 1 (headers) is
# untouched. Values below were captured from the source workbook and are
# written back to their rotated destination cell-by-cell. NumberFormat is
# forced to Text ("@") before writing text cells so that date/number-
# looking strings (e.g. "2023-08-22") are not auto-coerced by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- old row 2 -> new row 3 ----
$ws.Range("A3").Value = 111634177
$ws.Range("B3").Value = 90350
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = 'Ovaliderad'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 4786
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = 'Mandelriska'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = 'Lactarius volemus'
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = '(Fr.:Fr.) Fr.'
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = ''
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = ''
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = ''
$ws.Range("L3").ClearContents() | Out-Null
$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = ''
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = 'Hermansmåla, sydväst om Siggagölen, Bl'
$ws.Range("Q3").Value = 522930.7548289222
$ws.Range("R3").Value = 6247121.901725554
$ws.Range("S3").Value = 25
$ws.Range("T3").NumberFormat = "@"
$ws.Range("T3").Value = 'Blekinge'
$ws.Range("U3").NumberFormat = "@"
$ws.Range("U3").Value = 'Karlskrona'
$ws.Range("V3").NumberFormat = "@"
$ws.Range("V3").Value = 'Blekinge'
$ws.Range("W3").NumberFormat = "@"
$ws.Range("W3").Value = 'Tving'
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = '2023-08-22'
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = '00:00'
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = '2023-08-22'
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = '00:00'
$ws.Range("AC3").NumberFormat = "@"
$ws.Range("AC3").Value = 'Rikligt'
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AF3").NumberFormat = "@"
$ws.Range("AF3").Value = ''
$ws.Range("AG3").Value = $false
$ws.Range("AT3").NumberFormat = "@"
$ws.Range("AT3").Value = ''
$ws.Range("AW3").NumberFormat = "@"
$ws.Range("AW3").Value = 'Joakim Andersson Hemberg'
$ws.Range("AX3").NumberFormat = "@"
$ws.Range("AX3").Value = 'Joakim Andersson Hemberg'
$ws.Range("AY3").NumberFormat = "@"
$ws.Range("AY3").Value = ''

# ---- old row 3 -> new row 4 ----
$ws.Range("A4").Value = 111634171
$ws.Range("B4").Value = 73683
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = 'Ovaliderad'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 306
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = 'Kornig nållav'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = 'Chaenotheca chlorella'
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = '(Ach.) Müll.Arg.'
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = ''
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = ''
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = ''
$ws.Range("L4").ClearContents() | Out-Null
$ws.Range("N4").NumberFormat = "@"
$ws.Range("N4").Value = ''
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = 'Hermansmåla, söder om Siggagölen, Bl'
$ws.Range("Q4").Value = 522996.846862453
$ws.Range("R4").Value = 6247111.736777187
$ws.Range("S4").Value = 10
$ws.Range("T4").NumberFormat = "@"
$ws.Range("T4").Value = 'Blekinge'
$ws.Range("U4").NumberFormat = "@"
$ws.Range("U4").Value = 'Karlskrona'
$ws.Range("V4").NumberFormat = "@"
$ws.Range("V4").Value = 'Blekinge'
$ws.Range("W4").NumberFormat = "@"
$ws.Range("W4").Value = 'Tving'
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = '2023-08-22'
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = '00:00'
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = '2023-08-22'
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = '00:00'
$ws.Range("AC4").NumberFormat = "@"
$ws.Range("AC4").Value = 'På askhögstubbe.'
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AF4").NumberFormat = "@"
$ws.Range("AF4").Value = ''
$ws.Range("AG4").Value = $false
$ws.Range("AT4").NumberFormat = "@"
$ws.Range("AT4").Value = ''
$ws.Range("AW4").NumberFormat = "@"
$ws.Range("AW4").Value = 'Joakim Andersson Hemberg'
$ws.Range("AX4").NumberFormat = "@"
$ws.Range("AX4").Value = 'Joakim Andersson Hemberg'
$ws.Range("AY4").NumberFormat = "@"
$ws.Range("AY4").Value = ''

# ---- old row 4 -> new row 2 ----
$ws.Range("A2").Value = 111634202
$ws.Range("B2").Value = 92952
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = 'Ovaliderad'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = 'LC'
$ws.Range("E2").Value = 2779
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = 'Guldlockmossa'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = 'Homalothecium sericeum'
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = '(Hedw.) Schimp.'
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = ''
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = ''
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = ''
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = ''
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = ''
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = 'Hermansmåla, söder om Siggagölen, Bl'
$ws.Range("Q2").Value = 522996.846862453
$ws.Range("R2").Value = 6247111.736777187
$ws.Range("S2").Value = 10
$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value = 'Blekinge'
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = 'Karlskrona'
$ws.Range("V2").NumberFormat = "@"
$ws.Range("V2").Value = 'Blekinge'
$ws.Range("W2").NumberFormat = "@"
$ws.Range("W2").Value = 'Tving'
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = '2023-08-22'
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = '00:00'
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = '2023-08-22'
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = '00:00'
$ws.Range("AC2").NumberFormat = "@"
$ws.Range("AC2").Value = 'På ask.'
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AF2").NumberFormat = "@"
$ws.Range("AF2").Value = ''
$ws.Range("AG2").Value = $false
$ws.Range("AT2").NumberFormat = "@"
$ws.Range("AT2").Value = ''
$ws.Range("AW2").NumberFormat = "@"
$ws.Range("AW2").Value = 'Joakim Andersson Hemberg'
$ws.Range("AX2").NumberFormat = "@"
$ws.Range("AX2").Value = 'Joakim Andersson Hemberg'
$ws.Range("AY2").NumberFormat = "@"
$ws.Range("AY2").Value = ''
